# RoundRobinTesting.xlsx — "added comments and fixed it from RR-3 to RR-5"
#
# The second Round-Robin table (rows 8-13) is recomputed for a 5-job
# schedule instead of a 3-job one: the title changes from
# "Round Robin (3)" to "Round Robin (5)", the simulation numbers in rows
# 9-12 are updated, and the table is trimmed from column V back to
# column O (the old, wider 3-job layout isn't needed once the table only
# needs to track 5 jobs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title fix: Round Robin (3) -> Round Robin (5) ---------------------
$ws.Range("A8").Value = "Round Robin (5)"

# --- Recomputed simulation values --------------------------------------
$ws.Range("L9").Value  = 5
$ws.Range("M9").Value  = 2
$ws.Range("N9").Value  = 5

$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 13
$ws.Range("F10").Value = 5
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 7
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 8
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 2
$ws.Range("M10").Value = 3
$ws.Range("O10").Value = 0

$ws.Range("D11").Value = 5
$ws.Range("E11").Value = 10
$ws.Range("F11").Value = 15
$ws.Range("G11").Value = "Finished"
$ws.Range("H11").Value = 24
$ws.Range("I11").Value = "Finished"
$ws.Range("J11").Value = 31
$ws.Range("K11").Value = "Finished"
$ws.Range("L11").Value = 41
$ws.Range("M11").Value = 46
$ws.Range("O11").Value = "Finished"

$ws.Range("G12").Value = 19
$ws.Range("I12").Value = 26
$ws.Range("K12").Value = 36
$ws.Range("N12").Value = 48
$ws.Range("O12").Value = 51

# --- Table no longer extends past column O for rows 8-11; a few strays
#     in row 12 also get cleared -----------------------------------------
$ws.Range("P8:V11").ClearContents()
$ws.Range("L12").ClearContents()
$ws.Range("S12:T12").ClearContents()
$ws.Range("V12").ClearContents()

# --- Selection moves to A8 and the frozen/scrolled left column resets --
$ws.Range("A8").Select()
